# Change the table style (Table Design gallery pick) on the three tables
# that live on slides 14, 15 and 16 -- each slide's table is shape #1.
# {3E8A74F1-97AA-4589-AF02-19D64EEED89E} (the deck's custom "Table_0" style)
# -> {0E73C414-2940-4369-9B99-A0B0AF940F3F} (a built-in PowerPoint table style).
$p = $ppt.ActivePresentation

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    $tbl = $slide.Shapes.Item(1).Table
    $tbl.ApplyStyle("{0E73C414-2940-4369-9B99-A0B0AF940F3F}", $false)
}

# Switch the presentation's design from "Integral" (Red Violet colours) to
# the default "Office Theme" colours -- same effect as picking a different
# theme on the Design tab. Only the colour scheme actually differs between
# the two themes (fonts/format scheme are already identical), so drive the
# swap through the slide master's 12-slot colour scheme.
$master = $p.SlideMaster
$colorScheme = $master.ColorScheme

$officeColors = @(
    0,          # dk1
    16777215,   # lt1
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $colorScheme.Colors($i).RGB = $officeColors[$i - 1]
}
